$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 22.83000000000013
$ws.Range("G2").Value = [double]"3.881135413053016e-10"
$ws.Range("H2").Value = [double]"1.01155817323577e-09"
$ws.Range("K2").Value = 37.34033513073031
$ws.Range("L2").Value = "[23.35692331310119, 51.32374694835943]"
$ws.Range("M2").Value = [double]"4.801043305135977e-07"
$ws.Range("N2").Value = [double]"4.801043305135977e-07"
$ws.Range("O2").Value = 1.981184556317888
$ws.Range("P2").Value = "[1.578658170272349, 2.383710942363426]"
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 55.02280872851303
$ws.Range("T2").Value = "[47.68353551484772, 62.362081942178335]"
$ws.Range("W2").Value = 15.63135135135144
$ws.Range("X2").Value = 14.16876876876885
$ws.Range("Y2").Value = 17.09393393393403

# Row 3
$ws.Range("E3").Value = 23.55000000000024
$ws.Range("G3").Value = [double]"1.772437752123324e-11"
$ws.Range("H3").Value = [double]"2.075427131740883e-10"
$ws.Range("K3").Value = 40.87373995053681
$ws.Range("L3").Value = "[26.320638294365047, 55.42684160670857]"
$ws.Range("M3").Value = [double]"1.356191432488885e-07"
$ws.Range("N3").Value = [double]"2.712382864977769e-07"
$ws.Range("O3").Value = 2.345974093671657
$ws.Range("P3").Value = "[1.9937635058818093, 2.698184681461504]"
$ws.Range("S3").Value = 54.43193627748487
$ws.Range("T3").Value = "[46.97189737415021, 61.89197518081952]"
$ws.Range("W3").Value = 14.75705705705721
$ws.Range("X3").Value = 13.43693693693707
$ws.Range("Y3").Value = 16.07717717717735
